# issue #5: property land done
#
# Rewrites the 土地 (land) sheet's headers to the normalized English schema,
# appends the new property_category/category/date/legislator_*/source_file/index
# columns, cleans up stray full-width/half-width whitespace inside several
# shared text values on all three sheets, and normalizes a couple of
# thousand-separated numbers that were stored as (oddly punctuated) text.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

function Set-TextValue($ws, $row, $col, $text) {
    # Force the cell to stay a text cell (t="s") even when $text looks like
    # a pure number (e.g. "5707475") or an ISO date (e.g. "2012-11-28"), both
    # of which Excel would otherwise silently coerce to a number/date serial.
    # The classic leading-apostrophe text prefix avoids that, then the cell's
    # format is re-stamped from its left neighbour so the quote-prefix flag
    # doesn't leave a stray new cell style behind.
    $ws.Cells.Item($row, $col).Value = "'" + $text
    $ws.Cells.Item($row, $col - 1).Copy()
    $ws.Cells.Item($row, $col).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------
# Sheet "土地" (land)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("土地")

# Header row: translate the Chinese headers to the normalized English names.
$ws1.Cells.Item(1,2).Value = "name"
$ws1.Cells.Item(1,3).Value = "area"
$ws1.Cells.Item(1,4).Value = "share_portion"
$ws1.Cells.Item(1,5).Value = "owner"
$ws1.Cells.Item(1,6).Value = "register_date"
$ws1.Cells.Item(1,7).Value = "register_reason"
$ws1.Cells.Item(1,8).Value = "acquire_value"

# New trailing header columns I1:O1 (copy the bold/centered header style first).
$ws1.Cells.Item(1,8).Copy()
$ws1.Range("I1:O1").PasteSpecial($xlPasteFormats)
$ws1.Cells.Item(1,9).Value  = "property_category"
$ws1.Cells.Item(1,10).Value = "category"
$ws1.Cells.Item(1,11).Value = "date"
$ws1.Cells.Item(1,12).Value = "legislator_name"
$ws1.Cells.Item(1,13).Value = "legislator_id"
$ws1.Cells.Item(1,14).Value = "source_file"
$ws1.Cells.Item(1,15).Value = "index"

# Row 2 data cleanups + new trailing columns.
$ws1.Cells.Item(2,2).Value = "臺中市龍井區龍目井段水里社小段00080002地號"
$ws1.Cells.Item(2,6).Value = "83年11月16日"

$ws1.Cells.Item(2,8).Copy()
$ws1.Range("I2:O2").PasteSpecial($xlPasteFormats)
$ws1.Cells.Item(2,9).Value  = "land"
$ws1.Cells.Item(2,10).Value = "normal"
Set-TextValue $ws1 2 11 "2012-11-28"
$ws1.Cells.Item(2,12).Value = "顏清標"
$ws1.Cells.Item(2,13).Value = 979
$ws1.Cells.Item(2,14).Value = "tmp68961"
$ws1.Cells.Item(2,15).Value = 13

# Row 3 data cleanups + new trailing columns.
$ws1.Cells.Item(3,2).Value = "臺中市龍井區龍目井段水里社小段00080013±也號"
$ws1.Cells.Item(3,6).Value = "83年11月16日"

$ws1.Cells.Item(3,8).Copy()
$ws1.Range("I3:O3").PasteSpecial($xlPasteFormats)
$ws1.Cells.Item(3,9).Value  = "land"
$ws1.Cells.Item(3,10).Value = "normal"
Set-TextValue $ws1 3 11 "2012-11-28"
$ws1.Cells.Item(3,12).Value = "顏清標"
$ws1.Cells.Item(3,13).Value = 979
$ws1.Cells.Item(3,14).Value = "tmp68961"
$ws1.Cells.Item(3,15).Value = 14

# ---------------------------------------------------------------------
# Sheet "債務" (debt) -- text-only whitespace/number cleanups.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("債務")

$ws2.Cells.Item(2,4).Value  = "國泰世華銀行臺北市信義區松仁路"
Set-TextValue $ws2 2 5 "5707475"
$ws2.Cells.Item(2,6).Value  = "85年05月30日"

$ws2.Cells.Item(3,4).Value  = "華南銀行清水分行臺中市清水區中山路"
$ws2.Cells.Item(3,6).Value  = "84年07月08日"

$ws2.Cells.Item(4,4).Value  = "華南銀行清水分行臺中市清水區中山路"
$ws2.Cells.Item(4,6).Value  = "84年07月08日"
$ws2.Cells.Item(4,7).Value  = "繼承保證債務"

$ws2.Cells.Item(5,4).Value  = "元營建設股份有限公司臺中市沙鹿區北勢東路"
$ws2.Cells.Item(5,6).Value  = "97年10月23日"
$ws2.Cells.Item(5,7).Value  = "依據台灣台北地方法院97年10月23曰北院隆97執"

$ws2.Cells.Item(6,4).Value  = "華南銀行清水分行臺中市清水區中山路"
$ws2.Cells.Item(6,6).Value  = "84年07月08日"

$ws2.Cells.Item(7,4).Value  = "國泰世華銀行臺北市信義區松仁路"
$ws2.Cells.Item(7,6).Value  = "85年05月30日"

$ws2.Cells.Item(8,4).Value  = "國泰世華銀行臺北市信義區松仁路"
$ws2.Cells.Item(8,6).Value  = "85年05月30日"

$ws2.Cells.Item(9,4).Value  = "國泰世華銀行臺北市信義區松仁路"
$ws2.Cells.Item(9,6).Value  = "85年05月30日"

$ws2.Cells.Item(10,4).Value = "國泰世華銀行臺北市信義區松仁路"
$ws2.Cells.Item(10,6).Value = "85年11月16日"

# ---------------------------------------------------------------------
# Sheet "事業投資" (business investment) -- same style of cleanups.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("事業投資")

$ws3.Cells.Item(2,4).Value = "臺中市竹林里中山路紅竹巷58號1樓"
Set-TextValue $ws3 2 5 "1676000"
$ws3.Cells.Item(2,6).Value = "87年07月15日"

$ws3.Cells.Item(3,6).Value = "83年03月16日"

$ws3.Cells.Item(4,6).Value = "85年07月13日"
